$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 8226
$ws.Range("E2").Value = -187
$ws.Range("F2").Value = -187
$ws.Range("G2").Value = -101
$ws.Range("H2").Value = -74
$ws.Range("I2").Value = -68
$ws.Range("J2").Value = -6
$ws.Range("K2").Value = 7986
$ws.Range("L2").Value = 2583
$ws.Range("M2").Value = 5403
$ws.Range("N2").Value = 5402
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 913
$ws.Range("Q2").Value = 644
$ws.Range("R2").Value = -620
$ws.Range("S2").Value = 112
$ws.Range("T2").Value = 170
$ws.Range("U2").Value = 473
$ws.Range("V2").Value = 1042
$ws.Range("W2").Value = -2.28
$ws.Range("X2").Value = -0.9
$ws.Range("Y2").Value = -1.23
$ws.Range("Z2").Value = -0.93
$ws.Range("AA2").Value = 47.81
$ws.Range("AB2").Value = 496.68
$ws.Range("AC2").Value = -371
$ws.Range("AD2").Value = -78.95
$ws.Range("AE2").Value = 29826
$ws.Range("AF2").Value = 0.98
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 18252582

# Row 3
$ws.Range("D3").Value = 7929
$ws.Range("E3").Value = 422
$ws.Range("F3").Value = 422
$ws.Range("G3").Value = 428
$ws.Range("H3").Value = 349
$ws.Range("I3").Value = 349
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 9344
$ws.Range("L3").Value = 3661
$ws.Range("M3").Value = 5682
$ws.Range("N3").Value = 5681
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 913
$ws.Range("Q3").Value = 308
$ws.Range("R3").Value = -1000
$ws.Range("S3").Value = 551
$ws.Range("T3").Value = 251
$ws.Range("U3").Value = 57
$ws.Range("V3").Value = 1616
$ws.Range("W3").Value = 5.32
$ws.Range("X3").Value = 4.4
$ws.Range("Y3").Value = 6.29
$ws.Range("Z3").Value = 4.03
$ws.Range("AA3").Value = 64.44
$ws.Range("AB3").Value = 524.96
$ws.Range("AC3").Value = 1910
$ws.Range("AD3").Value = 18.32
$ws.Range("AE3").Value = 31435
$ws.Range("AF3").Value = 1.11
$ws.Range("AG3").Value = 750
$ws.Range("AH3").Value = 2.14
$ws.Range("AI3").Value = 38.87
$ws.Range("AJ3").Value = 18252582

# Row 4
$ws.Range("D4").Value = 8291
$ws.Range("E4").Value = -77
$ws.Range("F4").Value = -77
$ws.Range("G4").Value = -1
$ws.Range("H4").Value = -5
$ws.Range("I4").Value = -5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 8942
$ws.Range("L4").Value = 3534
$ws.Range("M4").Value = 5408
$ws.Range("N4").Value = 5407
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 913
$ws.Range("Q4").Value = -9
$ws.Range("R4").Value = -247
$ws.Range("S4").Value = 42
$ws.Range("T4").Value = 293
$ws.Range("U4").Value = -302
$ws.Range("V4").Value = 1895
$ws.Range("W4").Value = -0.92
$ws.Range("X4").Value = -0.06
$ws.Range("Y4").Value = -0.09
$ws.Range("Z4").Value = -0.06
$ws.Range("AA4").Value = 65.34999999999999
$ws.Range("AB4").Value = 507.72
$ws.Range("AC4").Value = -29
$ws.Range("AD4").Value = -915.5599999999999
$ws.Range("AE4").Value = 30565
$ws.Range("AF4").Value = 0.86
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 18252582

# Row 5
$ws.Range("D5").Value = 7732
$ws.Range("E5").Value = 187
$ws.Range("F5").Value = 187
$ws.Range("G5").Value = 203
$ws.Range("H5").Value = 159
$ws.Range("I5").Value = 150
$ws.Range("J5").Value = 9
$ws.Range("K5").Value = 10920
$ws.Range("L5").Value = 5290
$ws.Range("M5").Value = 5630
$ws.Range("N5").Value = 5554
$ws.Range("O5").Value = 76
$ws.Range("P5").Value = 913
$ws.Range("Q5").Value = 817
$ws.Range("R5").Value = -1158
$ws.Range("S5").Value = 398
$ws.Range("T5").Value = 238
$ws.Range("U5").Value = 578
$ws.Range("V5").Value = 2299
$ws.Range("W5").Value = 2.41
$ws.Range("X5").Value = 2.05
$ws.Range("Y5").Value = 2.74
$ws.Range("Z5").Value = 1.6
$ws.Range("AA5").Value = 93.97
$ws.Range("AB5").Value = 524.64
$ws.Range("AC5").Value = 824
$ws.Range("AD5").Value = 31.38
$ws.Range("AE5").Value = 31395
$ws.Range("AF5").Value = 0.82
$ws.Range("AG5").Value = 250
$ws.Range("AH5").Value = 0.97
$ws.Range("AI5").Value = 29.42
$ws.Range("AJ5").Value = 18252582

# Row 6
$ws.Range("D6").Value = 9140
$ws.Range("E6").Value = 51
$ws.Range("F6").Value = 51
$ws.Range("G6").Value = 79
$ws.Range("H6").Value = 58
$ws.Range("I6").Value = 49
$ws.Range("K6").Value = 9576
$ws.Range("L6").Value = 3892
$ws.Range("M6").Value = 5685
$ws.Range("N6").Value = 5600
$ws.Range("P6").Value = 913
$ws.Range("Q6").Value = 551
$ws.Range("R6").Value = -59
$ws.Range("S6").Value = -486
$ws.Range("T6").Value = 202
$ws.Range("U6").Value = 349
$ws.Range("V6").Value = 1861
$ws.Range("W6").Value = 0.5600000000000001
$ws.Range("X6").Value = 0.64
$ws.Range("Y6").Value = 0.87
$ws.Range("Z6").Value = 0.57
$ws.Range("AA6").Value = 68.45999999999999
$ws.Range("AB6").Value = 531.09
$ws.Range("AC6").Value = 267
$ws.Range("AD6").Value = 74.95999999999999
$ws.Range("AE6").Value = 31654
$ws.Range("AF6").Value = 0.63
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 18252582

# Row 7
$ws.Range("D7").Value = 8317
$ws.Range("E7").Value = -7
$ws.Range("G7").Value = 238
$ws.Range("H7").Value = 213
$ws.Range("I7").Value = 192
$ws.Range("K7").Value = 11982
$ws.Range("L7").Value = 5601
$ws.Range("M7").Value = 6380
$ws.Range("N7").Value = 5792
$ws.Range("P7").Value = 912
$ws.Range("Q7").Value = -771
$ws.Range("R7").Value = 325
$ws.Range("S7").Value = 1136
$ws.Range("T7").Value = 176
$ws.Range("U7").Value = -1825
$ws.Range("W7").Value = -0.08
$ws.Range("X7").Value = 2.56
$ws.Range("Y7").Value = 3.38
$ws.Range("Z7").Value = 1.98
$ws.Range("AA7").Value = 87.78
$ws.Range("AC7").Value = 1055
$ws.Range("AD7").Value = 19.58
$ws.Range("AE7").Value = 32742
$ws.Range("AF7").Value = 0.63
$ws.Range("AG7").Value = 120
$ws.Range("AH7").Value = 0.58
$ws.Range("AI7").Value = 11.38

# Row 8
$ws.Range("D8").Value = 9949
$ws.Range("E8").Value = 324
$ws.Range("G8").Value = 385
$ws.Range("H8").Value = 306
$ws.Range("I8").Value = 292
$ws.Range("K8").Value = 12330
$ws.Range("L8").Value = 5677
$ws.Range("M8").Value = 6654
$ws.Range("N8").Value = 6016
$ws.Range("P8").Value = 912
$ws.Range("Q8").Value = 319
$ws.Range("R8").Value = -318
$ws.Range("S8").Value = -425
$ws.Range("T8").Value = 215
$ws.Range("U8").Value = 245
$ws.Range("W8").Value = 3.26
$ws.Range("X8").Value = 3.07
$ws.Range("Y8").Value = 4.94
$ws.Range("Z8").Value = 2.51
$ws.Range("AA8").Value = 85.31999999999999
$ws.Range("AC8").Value = 1597
$ws.Range("AD8").Value = 13.71
$ws.Range("AE8").Value = 34005
$ws.Range("AF8").Value = 0.64
$ws.Range("AG8").Value = 290
$ws.Range("AH8").Value = 1.32
$ws.Range("AI8").Value = 18.16

# Row 9
$ws.Range("D9").Value = 9494
$ws.Range("E9").Value = 353
$ws.Range("G9").Value = 418
$ws.Range("H9").Value = 327
$ws.Range("I9").Value = 314
$ws.Range("K9").Value = 12454
$ws.Range("L9").Value = 5517
$ws.Range("M9").Value = 6937
$ws.Range("N9").Value = 6299
$ws.Range("P9").Value = 912
$ws.Range("Q9").Value = 494
$ws.Range("R9").Value = -312
$ws.Range("S9").Value = -255
$ws.Range("T9").Value = 210
$ws.Range("U9").Value = 185
$ws.Range("W9").Value = 3.72
$ws.Range("X9").Value = 3.44
$ws.Range("Y9").Value = 5.11
$ws.Range("Z9").Value = 2.64
$ws.Range("AA9").Value = 79.53
$ws.Range("AC9").Value = 1723
$ws.Range("AD9").Value = 12.71
$ws.Range("AE9").Value = 35608
$ws.Range("AF9").Value = 0.62
$ws.Range("AG9").Value = 350
$ws.Range("AH9").Value = 1.6
$ws.Range("AI9").Value = 20.31
